# PagoAnticipadoPlazo.xlsx - "update entregable 1, 2"
#
# Row1 headers D:F get new text ("Estado"/"Transaccion"/"Fecha") while keeping
# the orange header fill (same style already used by F1/G1/H1); row2's test
# result columns (D:H) are cleared out and A2's credential value is replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers -----------------------------------------------------
# D1 needs the "Estado" header with the orange-fill style that F1 already
# has -> copy F1's format into D1, then overwrite the text.
$ws.Range("F1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "Estado"

# E1 needs the "Transaccion" header, same orange-fill style as G1.
$ws.Range("G1").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Transaccion"

# F1 keeps its own (already-orange) style, just gets new text.
$ws.Range("F1").Value = "Fecha"

# Old G1/H1 headers ("Transaccion"/"Fecha") are no longer needed.
$ws.Range("G1:H1").Clear()

# --- Row 2 data ----------------------------------------------------------
$ws.Range("A2").Value = "ebenito"

# The old run's result columns (D2:H2) are removed entirely.
$ws.Range("D2:H2").Clear()

# --- Cosmetic bits --------------------------------------------------------
# Columns shrink/grow to fit their new (shorter/longer) header text.
$ws.Columns("A").ColumnWidth = 6.25
$ws.Columns("D").ColumnWidth = 6.25
$ws.Columns("E").ColumnWidth = 17

$ws.Range("I6").Select()
